$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.2307692307692308
$ws.Cells.Item(2, 3).Value = 0.4585798816568047
$ws.Cells.Item(2, 10).Value = 0.02958579881656805
$ws.Cells.Item(2, 16).Value = 0.1597633136094675
$ws.Cells.Item(2, 19).Value = 0.121301775147929
$ws.Cells.Item(3, 2).Value = 0.006097560975609756
$ws.Cells.Item(3, 3).Value = 0.03658536585365853
$ws.Cells.Item(3, 10).Value = 0.06097560975609756
$ws.Cells.Item(3, 16).Value = 0.6890243902439024
$ws.Cells.Item(3, 19).Value = 0.2073170731707317
$ws.Cells.Item(4, 10).Value = 0.08333333333333333
$ws.Cells.Item(4, 16).Value = 0.6875
$ws.Cells.Item(4, 19).Value = 0.2291666666666667
$ws.Cells.Item(6, 2).Value = 0.0365296803652968
$ws.Cells.Item(6, 4).Value = 0.0136986301369863
$ws.Cells.Item(6, 6).Value = 0.0867579908675799
$ws.Cells.Item(6, 10).Value = 0.2602739726027397
$ws.Cells.Item(6, 15).Value = 0.0228310502283105
$ws.Cells.Item(6, 17).Value = 0.1643835616438356
$ws.Cells.Item(6, 18).Value = 0.0730593607305936
$ws.Cells.Item(6, 19).Value = 0.3424657534246575
$ws.Cells.Item(7, 2).Value = 0.1085714285714286
$ws.Cells.Item(7, 4).Value = 0.01142857142857143
$ws.Cells.Item(7, 6).Value = 0.01142857142857143
$ws.Cells.Item(7, 10).Value = 0.1257142857142857
$ws.Cells.Item(7, 15).Value = 0.03428571428571429
$ws.Cells.Item(7, 17).Value = 0.2114285714285714
$ws.Cells.Item(7, 18).Value = 0.1085714285714286
$ws.Cells.Item(7, 19).Value = 0.3885714285714286
$ws.Cells.Item(8, 2).Value = 0.08858267716535433
$ws.Cells.Item(8, 4).Value = 0.00984251968503937
$ws.Cells.Item(8, 5).Value = 0.001968503937007874
$ws.Cells.Item(8, 6).Value = 0.05511811023622047
$ws.Cells.Item(8, 10).Value = 0.1299212598425197
$ws.Cells.Item(8, 15).Value = 0.01377952755905512
$ws.Cells.Item(8, 17).Value = 0.1909448818897638
$ws.Cells.Item(8, 18).Value = 0.1023622047244094
$ws.Cells.Item(8, 19).Value = 0.4074803149606299
$ws.Cells.Item(9, 2).Value = 0.09547738693467336
$ws.Cells.Item(9, 4).Value = 0.01507537688442211
$ws.Cells.Item(9, 6).Value = 0.06532663316582915
$ws.Cells.Item(9, 10).Value = 0.1206030150753769
$ws.Cells.Item(9, 15).Value = 0.02512562814070352
$ws.Cells.Item(9, 17).Value = 0.2060301507537688
$ws.Cells.Item(9, 18).Value = 0.1055276381909548
$ws.Cells.Item(9, 19).Value = 0.3668341708542713
$ws.Cells.Item(10, 2).Value = 0.09661538461538462
$ws.Cells.Item(10, 4).Value = 0.02215384615384615
$ws.Cells.Item(10, 5).Value = 0.0006153846153846154
$ws.Cells.Item(10, 6).Value = 0.04615384615384616
$ws.Cells.Item(10, 10).Value = 0.1212307692307692
$ws.Cells.Item(10, 15).Value = 0.01292307692307692
$ws.Cells.Item(10, 17).Value = 0.2258461538461538
$ws.Cells.Item(10, 18).Value = 0.1009230769230769
$ws.Cells.Item(10, 19).Value = 0.3735384615384615
$ws.Cells.Item(11, 7).Value = 0.1508771929824561
$ws.Cells.Item(11, 10).Value = 0.09824561403508772
$ws.Cells.Item(11, 11).Value = 0.2210526315789474
$ws.Cells.Item(11, 12).Value = 0.5228070175438596
$ws.Cells.Item(11, 19).Value = 0.007017543859649123
$ws.Cells.Item(12, 7).Value = 0.7702702702702703
$ws.Cells.Item(12, 10).Value = 0.1756756756756757
$ws.Cells.Item(12, 11).Value = 0.006756756756756757
$ws.Cells.Item(12, 12).Value = 0.01351351351351351
$ws.Cells.Item(12, 19).Value = 0.03378378378378379
$ws.Cells.Item(13, 7).Value = 0.6285714285714286
$ws.Cells.Item(13, 10).Value = 0.2857142857142857
$ws.Cells.Item(13, 19).Value = 0.08571428571428572
$ws.Cells.Item(14, 7).Value = 0.5
$ws.Cells.Item(14, 10).Value = 0.5
$ws.Cells.Item(15, 6).Value = 0.0211864406779661
$ws.Cells.Item(15, 8).Value = 0.1483050847457627
$ws.Cells.Item(15, 9).Value = 0.05508474576271186
$ws.Cells.Item(15, 10).Value = 0.3983050847457627
$ws.Cells.Item(15, 11).Value = 0.03389830508474576
$ws.Cells.Item(15, 13).Value = 0.01271186440677966
$ws.Cells.Item(15, 15).Value = 0.05508474576271186
$ws.Cells.Item(15, 19).Value = 0.2754237288135593
$ws.Cells.Item(16, 6).Value = 0.01030927835051546
$ws.Cells.Item(16, 8).Value = 0.1597938144329897
$ws.Cells.Item(16, 9).Value = 0.04123711340206185
$ws.Cells.Item(16, 10).Value = 0.5154639175257731
$ws.Cells.Item(16, 11).Value = 0.08247422680412371
$ws.Cells.Item(16, 13).Value = 0.02577319587628866
$ws.Cells.Item(16, 14).Value = 0.005154639175257732
$ws.Cells.Item(16, 15).Value = 0.04123711340206185
$ws.Cells.Item(16, 19).Value = 0.1185567010309278
$ws.Cells.Item(17, 6).Value = 0.02292768959435626
$ws.Cells.Item(17, 8).Value = 0.1798941798941799
$ws.Cells.Item(17, 9).Value = 0.1111111111111111
$ws.Cells.Item(17, 10).Value = 0.4673721340388007
$ws.Cells.Item(17, 11).Value = 0.06349206349206349
$ws.Cells.Item(17, 13).Value = 0.007054673721340388
$ws.Cells.Item(17, 15).Value = 0.04761904761904762
$ws.Cells.Item(17, 19).Value = 0.1005291005291005
$ws.Cells.Item(18, 6).Value = 0.01107011070110701
$ws.Cells.Item(18, 8).Value = 0.1881918819188192
$ws.Cells.Item(18, 9).Value = 0.07749077490774908
$ws.Cells.Item(18, 10).Value = 0.4833948339483395
$ws.Cells.Item(18, 11).Value = 0.08118081180811808
$ws.Cells.Item(18, 13).Value = 0.01107011070110701
$ws.Cells.Item(18, 14).Value = 0.007380073800738007
$ws.Cells.Item(18, 15).Value = 0.06273062730627306
$ws.Cells.Item(18, 19).Value = 0.07749077490774908
$ws.Cells.Item(19, 6).Value = 0.02121640735502122
$ws.Cells.Item(19, 8).Value = 0.2072135785007072
$ws.Cells.Item(19, 9).Value = 0.06577086280056577
$ws.Cells.Item(19, 10).Value = 0.417963224893918
$ws.Cells.Item(19, 11).Value = 0.09900990099009901
$ws.Cells.Item(19, 13).Value = 0.01768033946251768
$ws.Cells.Item(19, 14).Value = 0.0007072135785007072
$ws.Cells.Item(19, 15).Value = 0.07001414427157002
$ws.Cells.Item(19, 19).Value = 0.1004243281471004
